# "add USCDI Goals and Preferences part 1"
# Adds two new profile rows (Treatment Intervention Preference, Care
# Experience Preference) near the top of the "profiles" sheet, and makes
# "profiles" the active tab/sheet.

$wb = $excel.ActiveWorkbook

# --- Bump the internal sheetId of "profiles" (33 -> 34) -------------------
# Real Excel re-created this sheet (sheetId goes from 33 to 34 in the
# target file) while everything else about it (position, formatting,
# data) stays the same. Copying the sheet right next to itself clones all
# of its content/formatting under a fresh sheetId, then we delete the
# original and rename the clone back to "profiles" so it ends up in the
# same tab position it started in.
$oldProfiles = $wb.Worksheets.Item("profiles")
$oldProfiles.Copy($null, $oldProfiles)
$wb.Worksheets.Item("profiles").Delete()
$clone = $wb.Worksheets.Item("profiles (2)")
$clone.Name = "profiles"

# --- Insert the two new profile rows --------------------------------------
$ws = $wb.Worksheets.Item("profiles")
$ws.Rows("32:33").Insert()

$ws.Range("A32").Value = "http://hl7.org/fhir/us/core/StructureDefinition/us-core-treatment-intervention-preference"
$ws.Range("B32").Value = "US Core Treatment Intervention Preference Profile"
$ws.Range("D32").Value = "SHALL"
$ws.Range("E32").Value = "Observation"

$ws.Range("A33").Value = "http://hl7.org/fhir/us/core/StructureDefinition/us-core-care-experience-preference"
$ws.Range("B33").Value = "US Core Care Experience Preference Profile"
$ws.Range("D33").Value = "SHALL"
$ws.Range("E33").Value = "Observation"

# --- Update the view state: "profiles" becomes the active/selected tab ----
$ws.Activate()
$ws.Range("A44").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
